# 3DES Adicionada aula 1 RMS
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FREQ")

# Fill in attendance ("P" = present, "F" = falta/absent) for the new
# "aula 1" column (G) added to the RMST/PROJ block.
$ws.Range("G3").Value = "P"
$ws.Range("G4").Value = "P"
$ws.Range("G5").Value = "P"
$ws.Range("G6").Value = "P"
$ws.Range("G7").Value = "F"
$ws.Range("G8").Value = "P"
$ws.Range("G9").Value = "P"
$ws.Range("G10").Value = "F"
$ws.Range("G11").Value = "P"
$ws.Range("G12").Value = "F"
$ws.Range("G13").Value = "P"
$ws.Range("G14").Value = "P"
$ws.Range("G15").Value = "F"
$ws.Range("G16").Value = "P"
$ws.Range("G17").Value = "P"
$ws.Range("G18").Value = "P"
$ws.Range("G19").Value = "P"
$ws.Range("G20").Value = "P"

# Reset the active cell back to A1 so the sheet no longer persists a
# stale selection at J4.
$ws.Activate()
$ws.Range("A1").Select()
